{"js": "// Replace the multiplication expressions in the practice-sheet table.\n// Each cell holds a single unique \"AAA\u00d7B=\" run, so an exact text\n// search-and-replace for each old/new pair is safe and unambiguous.\nconst replacements = [\n  [\"400\u00d75=\", \"392\u00d76=\"],\n  [\"882\u00d72=\", \"176\u00d74=\"],\n  [\"452\u00d77=\", \"851\u00d76=\"],\n  [\"311\u00d79=\", \"619\u00d78=\"],\n  [\"355\u00d74=\", \"469\u00d77=\"],\n  [\"256\u00d72=\", \"209\u00d77=\"],\n  [\"932\u00d76=\", \"885\u00d77=\"],\n  [\"325\u00d77=\", \"547\u00d77=\"],\n  [\"792\u00d77=\", \"917\u00d72=\"],\n  [\"190\u00d76=\", \"264\u00d78=\"],\n  [\"452\u00d75=\", \"680\u00d73=\"],\n  [\"532\u00d78=\", \"361\u00d78=\"],\n  [\"196\u00d73=\", \"809\u00d74=\"],\n  [\"157\u00d76=\", \"494\u00d73=\"],\n  [\"533\u00d75=\", \"259\u00d78=\"],\n  [\"218\u00d79=\", \"504\u00d79=\"],\n  [\"799\u00d72=\", \"591\u00d76=\"],\n  [\"982\u00d73=\", \"352\u00d76=\"],\n  [\"274\u00d72=\", \"185\u00d75=\"],\n  [\"717\u00d76=\", \"875\u00d72=\"],\n  [\"119\u00d76=\", \"396\u00d78=\"],\n  [\"223\u00d75=\", \"897\u00d76=\"],\n  [\"633\u00d79=\", \"550\u00d75=\"],\n  [\"264\u00d74=\", \"501\u00d73=\"],\n  [\"477\u00d74=\", \"864\u00d74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication expressions in the practice-sheet table.\n# Each cell holds a single unique \"AAA\u00d7B=\" run, so an exact text\n# find & replace for each old/new pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"400\u00d75=\", \"392\u00d76=\"),\n    @(\"882\u00d72=\", \"176\u00d74=\"),\n    @(\"452\u00d77=\", \"851\u00d76=\"),\n    @(\"311\u00d79=\", \"619\u00d78=\"),\n    @(\"355\u00d74=\", \"469\u00d77=\"),\n    @(\"256\u00d72=\", \"209\u00d77=\"),\n    @(\"932\u00d76=\", \"885\u00d77=\"),\n    @(\"325\u00d77=\", \"547\u00d77=\"),\n    @(\"792\u00d77=\", \"917\u00d72=\"),\n    @(\"190\u00d76=\", \"264\u00d78=\"),\n    @(\"452\u00d75=\", \"680\u00d73=\"),\n    @(\"532\u00d78=\", \"361\u00d78=\"),\n    @(\"196\u00d73=\", \"809\u00d74=\"),\n    @(\"157\u00d76=\", \"494\u00d73=\"),\n    @(\"533\u00d75=\", \"259\u00d78=\"),\n    @(\"218\u00d79=\", \"504\u00d79=\"),\n    @(\"799\u00d72=\", \"591\u00d76=\"),\n    @(\"982\u00d73=\", \"352\u00d76=\"),\n    @(\"274\u00d72=\", \"185\u00d75=\"),\n    @(\"717\u00d76=\", \"875\u00d72=\"),\n    @(\"119\u00d76=\", \"396\u00d78=\"),\n    @(\"223\u00d75=\", \"897\u00d76=\"),\n    @(\"633\u00d79=\", \"550\u00d75=\"),\n    @(\"264\u00d74=\", \"501\u00d73=\"),\n    @(\"477\u00d74=\", \"864\u00d74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
